# Collections Ch 7 & 8: fixed numbers; tables; charts; deleted old charts
#
# This script reproduces the content-level edits on slide 1 of the
# "product-arraylist" deck: two text corrections (the stray duplicate
# count label, and the rounded memory-size figure) plus the small
# shape repositioning/resizing that PowerPoint performs when the file
# is touched/resaved (values snap to 1/8-point increments).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------
# 1) Regular Pentagon 6 (Id 7)
# ---------------------------------------------------------------
$shp = $s.Shapes.Item(1)
$shp.Left   = 261.75
$shp.Top    = 203.87503937007875
$shp.Width  = 178.12503937007875
$shp.Height = 67.25

# ---------------------------------------------------------------
# 2) Rectangle 7 (Id 8)
# ---------------------------------------------------------------
$shp = $s.Shapes.Item(2)
$shp.Left   = 285.25
$shp.Top    = 343.5
$shp.Width  = 140.5
$shp.Height = 46.75

# ---------------------------------------------------------------
# 3) Straight Arrow Connector 11 (Id 12)
# ---------------------------------------------------------------
$shp = $s.Shapes.Item(3)
$shp.Left   = 322.5
$shp.Top    = 175.62503937007875
$shp.Width  = 56.62496062992126
$shp.Height = 0.12503937007874016

# ---------------------------------------------------------------
# 4) Straight Arrow Connector 13 (Id 14)
# ---------------------------------------------------------------
$shp = $s.Shapes.Item(4)
$shp.Left   = 313.75
$shp.Top    = 308.12503937007875
$shp.Width  = 74.12496062992126
$shp.Height = 0.12503937007874016

# ---------------------------------------------------------------
# 5) TextBox 27 (Id 28) - "Supplier / x400,000" label
#    text unchanged, font normalized to Calibri
# ---------------------------------------------------------------
$shp = $s.Shapes.Item(5)
$shp.Left   = 293.0
$shp.Top    = 345.25
$shp.Width  = 115.62503937007874
$shp.Height = 41.25
$shp.TextFrame.TextRange.Font.Name = "Calibri"

# ---------------------------------------------------------------
# 6) TextBox 35 (Id 36) - count label corrected "100,000" -> "1"
# ---------------------------------------------------------------
$shp = $s.Shapes.Item(6)
$shp.Left   = 350.87503937007875
$shp.Top    = 157.0
$shp.Width  = 21.62496062992126
$shp.Height = 24.0
$shp.TextFrame.TextRange.Text = ""
$shp.TextFrame.TextRange.Text = "1"
$shp.TextFrame.TextRange.Font.Name = "Calibri"

# ---------------------------------------------------------------
# 7) TextBox 36 (Id 37) - "4" label, text unchanged, font normalized
# ---------------------------------------------------------------
$shp = $s.Shapes.Item(7)
$shp.Left   = 350.75
$shp.Top    = 300.25
$shp.Width  = 21.75
$shp.Height = 24.12503937007874
$shp.TextFrame.TextRange.Font.Name = "Calibri"

# ---------------------------------------------------------------
# 8) TextBox 32 (Id 33) - "ArrayList / x100,000 = 7.629MB"
#    fixed the rounded memory size to "7.6MB" and simplified runs
# ---------------------------------------------------------------
$shp = $s.Shapes.Item(8)
$shp.Left   = 282.0
$shp.Top    = 221.87503937007875
$shp.Width  = 117.5
$shp.Height = 40.75
$shp.TextFrame.TextRange.Text = ""
$shp.TextFrame.TextRange.Text = "         ArrayList`rx100,000 = 7.6MB"
$tr = $shp.TextFrame.TextRange
$tr.Font.Name = "Calibri"
$secondLine = $tr.Characters(20, 16)
$secondLine.Font.Italic = $true
$secondLine.Font.Name = "Calibri"

# ---------------------------------------------------------------
# 9) TextBox 46 (Id 47) - "Product / x100,000" label, text unchanged
# ---------------------------------------------------------------
$shp = $s.Shapes.Item(9)
$shp.Left   = 285.25
$shp.Top    = 100.62503937007874
$shp.Width  = 94.75
$shp.Height = 41.12496062992126
$shp.TextFrame.TextRange.Font.Name = "Calibri"

# ---------------------------------------------------------------
# 10) Rectangle 48 (Id 49)
# ---------------------------------------------------------------
$shp = $s.Shapes.Item(10)
$shp.Left   = 273.37503937007875
$shp.Top    = 100.62503937007874
$shp.Width  = 144.62496062992125
$shp.Height = 46.75
